$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows in column C that currently hold the shared string "space" and
# should instead hold the numeric value 1 (bug fix for N-back corrAns).
$rows = @(7,12,15,22,25,31,34,40,43,50,53,58,67,72,75,82,85,91,94,100,103,110,113,118)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = 1
}

# Restore the active sheet selection to J31, as captured in the saved view.
$ws.Activate()
$ws.Range("J31").Select()
